# Lesson_3.pptx edit: "Add lesson4 and exercise4"
#
# 1) Slide 11 ("Let's Start to Play") - the italic sub-title line is
#    changed from "Programming a virtual OS in a Structured Language" to
#    "Programming a Graphical Application of the Pomodoro Timer", split
#    across three runs (the new wording now references the Pomodoro app
#    that exercise 4 is built around).
# 2) Slide 13 ("Exercise 3 part 2") - two pairs of runs that were
#    previously split are merged back into single runs:
#      "Pomodoro " + "Config"                       -> "Pomodoro Config"
#      "./" + "ComputerProgrammingBasic/exercise3/config.py"
#                                                      -> "./ComputerProgrammingBasic/exercise3/config.py"

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Slide 11 title line rewrite
# ---------------------------------------------------------------------
$slide11 = $p.Slides.Item(11)

$titleShape = $null
for ($i = 1; $i -le $slide11.Shapes.Count; $i++) {
    $shape = $slide11.Shapes.Item($i)
    if ($shape.HasTextFrame) {
        if ($shape.TextFrame.TextRange.Text -like "*Programming a virtual OS*") {
            $titleShape = $shape
        }
    }
}

$titleRange = $titleShape.TextFrame.TextRange
$titleRange.Text = "Programming a Graphical Application of the Pomodoro Timer"

$totalLen = $titleRange.Text.Length
# "Programming " (12 chars) | "a " (2 chars) | "Graphical Application of the Pomodoro Timer" (rest)
$runB = $titleRange.Characters(13, 2)
$runB.Font.Italic = -1
$runC = $titleRange.Characters(15, $totalLen - 14)
$runC.Font.Italic = -1

# ---------------------------------------------------------------------
# 2) Slide 13 run merges
# ---------------------------------------------------------------------
$slide13 = $p.Slides.Item(13)

$bodyShape = $null
for ($i = 1; $i -le $slide13.Shapes.Count; $i++) {
    $shape = $slide13.Shapes.Item($i)
    if ($shape.HasTextFrame) {
        if ($shape.TextFrame.TextRange.Text -like "*Pomodoro*config.py*") {
            $bodyShape = $shape
        }
    }
}

$bodyRange = $bodyShape.TextFrame.TextRange

# "Review the Pomodoro Config." -> merge "Pomodoro " + "Config" runs
$pomodoroConfig = $bodyRange.Characters(12, 15)
$pomodoroConfig.Text = "Pomodoro Config"

# "./ComputerProgrammingBasic/exercise3/config.py" -> merge "./" + rest
$configPath = $bodyRange.Characters(28, 46)
$configPath.Text = "./ComputerProgrammingBasic/exercise3/config.py"
